$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 55
$ws.Range("E12").Value = "55/140"
